$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: add End Hour/Minute (finish the existing lecture 12 finger exercises entry)
$ws.Range("D12").Value = 22
$ws.Range("E12").Value = 38

# Row 13: new entry for a second "CS introduction Lecture 12" session (typo'd task name)
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = 45813
$ws.Range("B13").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = 15
$ws.Range("E13").Value = 47
$ws.Range("F13").Value = "CS intoduction Lecture 12"

# Update selection to match resulting document state
$ws.Range("A14").Select()
